$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the formatting (borders/fill) used by the existing data rows (3-7)
# down into the three new rows before filling in their content.
$ws.Range("A7:B7").Copy() | Out-Null
$ws.Range("A8:B10").PasteSpecial(-4122) | Out-Null

# Row 8: Sl.No 6 - Write data pin
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "To set and check the Write data pin when port direction is set as input."
$ws.Range("B8").Characters(21, 15).Font.Bold = $true

# Row 9: Sl.No 7 - Read data pin
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "To set and check the Read data pin when port direction is set as output."
$ws.Range("B9").Characters(21, 14).Font.Bold = $true

# Row 10: Sl.No 8 - pull up or pull down configuration
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "To enable or disable the pull up or pull down configuration"
$ws.Range("B10").Characters(26, 20).Font.Bold = $true

# Match the final selection left behind by the author (scrolled past the new rows)
$ws.Range("B15").Select() | Out-Null
